$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Mark rows 50-54 (Task ID 49-53) as Completed = Yes in column E.
#    Copy the formatting used by the existing "Yes" cells (e.g. E49) so the
#    new cells pick up the same style index, then set the value.
$ws.Range("E49").Copy()
$ws.Range("E50:E54").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E50").Value = "Yes"
$ws.Range("E51").Value = "Yes"
$ws.Range("E52").Value = "Yes"
$ws.Range("E53").Value = "Yes"
$ws.Range("E54").Value = "Yes"

# 2) Append 5 new rows (Josh's tasks) to Table_1, which extends the table
#    range and the sheet data accordingly.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Copy the date format from the last pre-existing data row (D54) onto the
# new date cells so they keep the same number format / style.
$ws.Range("D54").Copy()
$ws.Range("D55:D59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Task names (column B) are entered in the same order the author originally
# typed them so new shared-string entries land at the same indices.
$ws.Range("B55").Value = "Finish up the poster"
$ws.Range("B56").Value = "Demonstrate LED driver on Human Brain"
$ws.Range("B59").Value = "Work on Final report and any outstanding documentation"
$ws.Range("B57").Value = "Work on LED driver on Raspberry pi"
$ws.Range("B58").Value = "Implement Josh's software onto the Raspberry Pi"

$ws.Range("A55").Value = 54
$ws.Range("C55").Value = "George Proios"
$ws.Range("D55").Value = 43384

$ws.Range("A56").Value = 55
$ws.Range("C56").Value = "Josh Francis"
$ws.Range("D56").Value = 43384

$ws.Range("A57").Value = 56
$ws.Range("C57").Value = "Yeqing Liu"
$ws.Range("D57").Value = 43384

$ws.Range("A58").Value = 57
$ws.Range("C58").Value = "Michael Douglas"
$ws.Range("D58").Value = 43384

$ws.Range("A59").Value = 58
$ws.Range("C59").Value = "Timothy Finn"
$ws.Range("D59").Value = 43384

# 3) Update the current selection to match where the user ended up editing.
$ws.Range("A55:D59").Select()
